$d = $word.ActiveDocument

# 1. Update the "Software Development and Innovation" heading text
$d.Content.Find.Execute("Software Development and Innovation", $true, $false, $false, $false, $false, $true, 1, $false, "Technical Innovation & Platform Development", 2) | Out-Null

# 2. Expand the first bullet and insert four new bullets after it
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Conceived and deployed redistricting software used by thousands of analysts nationwide")) {
        $target = $p
        break
    }
}
$target.Range.Text = "• Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide"
$newBullets = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
$newBullets += "`r• Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors"
$newBullets += "`r• Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls"
$newBullets += "`r• Created SimCrisis platform for humanitarian intervention modeling used by International Red Cross and UNICEF"
$target.Range.InsertAfter("`r" + $newBullets)

# 3. Replace the "Created econometric simulation..." bullet with the new
#    "Data Engineering & Analytics" heading, then insert its bullets plus
#    the "Research Leadership & Client Success" heading and its bullets.
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Created econometric simulation platform for humanitarian intervention modeling")) {
        $target2 = $p
        break
    }
}
$target2.Range.Text = "Data Engineering & Analytics"

$deBullets = "• Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change"
$deBullets += "`r• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%"
$deBullets += "`r• Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling"
$deBullets += "`r• Built fraud detection systems for campaign finance data analysis across multi-terabyte datasets"
$deBullets += "`r• Transformed small data team into big data engineering team using Hadoop Clusters and Hive on AWS"
$deBullets += "`r• Introduced version control and Agile methodologies, improving project delivery timelines by 40%"
$deBullets += "`rResearch Leadership & Client Success"
$deBullets += "`r• Led multi-million dollar research projects involving sensitive consumer data with privacy compliance"
$deBullets += "`r• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders"
$deBullets += "`r• Delivered actionable consumer insights and market intelligence for political candidates and major organizations"
$target2.Range.InsertAfter("`r" + $deBullets)

# Apply Heading 3 style to the two new heading paragraphs (after all text
# insertion is done, so the style does not bleed into sibling paragraphs).
$target2.Style = "Heading 3"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Research Leadership & Client Success")) {
        $p.Style = "Heading 3"
        break
    }
}

# 4. After the "Built comprehensive survey operations..." bullet, append two
#    more bullets about expert testimony.
$target3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Built comprehensive survey operations platform from RFP through deployment")) {
        $target3 = $p
        break
    }
}
$testimonyBullets = "• Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership"
$testimonyBullets += "`r• Redistricting analysis used in court cases with rigorous methodology and expert testimony"
$target3.Range.InsertAfter("`r" + $testimonyBullets)
